$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component Bank")

# --- Update the "Components from Component Bank" table with the new
# --- requirements (rows 8-13). Written in the same order the shared
# --- strings ended up being introduced so text/order matches.

# Row 13 (will become "Logitech c920 HD webcam") - currently empty/unused
$ws.Range("A13").Value = "Logitech c920 HD webcam"

# Row 8: "Stepper Motors" (x3 required) -> specific Wantai NEMA8 stepper, qty 1
$ws.Range("A8").Value = "Wantai Stepper Motor 20BYGH406, NEMA8, 1.8 deg/step"
$ws.Range("B8").Value = 1
$ws.Range("C8").Style = "Normal"
$ws.Range("C8").Value = "Required"
$ws.Range("C8").HorizontalAlignment = -4131
$ws.Range("C8").VerticalAlignment = -4108

# Row 10: "Servo Motor" -> specific Wantai dual shaft NEMA17 stepper, qty 1
$ws.Range("A10").Value = "Wantai Dual Shaft Stepper Motor 42BYGHW920L21B2, NEMA17, 1.8 deg/step"

# Row 9: "Stepper Motor Drivers" (x3 required) -> specific Wantai NEMA14 stepper, qty 1
$ws.Range("A9").Value = "Wantai Stepper Motor 35BYG312P1,NEMA14, 1.8 deg/step"
$ws.Range("B9").Value = 1

# Row 12 (new row): DS3115MG Servo, qty 1, Required
$ws.Range("A12").Value = "DS3115MG Servo"

# Row 11: "Machine Vision Camera" -> Pololu stepper driver carrier, qty 3
$ws.Range("A11").Value = "Pololu DRV8825 Stepper Motor Driver Carrier, High Current"
$ws.Range("B11").Value = 3

# Fill remaining Quantity / Status cells for rows 9-13
$ws.Range("C9").Value = "Required"
$ws.Range("C10").Value = "Required"
$ws.Range("C11").Value = "Required"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Required"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Required"

# Match formatting of the rest of the table for the two brand-new rows
$fmtRange = $ws.Range("A12:C13")
$fmtRange.HorizontalAlignment = -4131
$fmtRange.VerticalAlignment = -4108

# Expand the bound table (Table3) to cover the two new rows
$t = $ws.ListObjects.Item("Table3")
$t.Resize($ws.Range("A4:C13"))

# Column A widened automatically by Excel to fit the longer descriptions
$ws.Columns.Item(1).ColumnWidth = 66.25

# Leave the cursor where the author last left it
$ws.Range("F16").Select()
